# [JH] 13. UI&Data Work
# SwordWoman.xlsx: the "Class" column (N) for every character-level row
# (rows 5-44) no longer stores the literal text "Knight" - it now holds
# the numeric placeholder 0. Because "Knight" becomes unused once every
# reference to it is gone, it naturally drops out of the shared string
# table, shifting the remaining string indices down - which is exactly
# what the target workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 5; $row -le 44; $row++) {
    $ws.Cells.Item($row, 14).Value = 0
}

# Restore the cursor/selection to where the author last left it.
$ws.Range("C5").Select() | Out-Null

Write-Host "Updated Class column (N5:N44) to 0 and reselected C5"
